$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.30261196328375123,
    0.22278057385692498,
    0.68958690916092702,
    0.50418966227800732,
    0.35789267429157523,
    0.45082170644060554,
    0.50582692017245301,
    0.75990210820240955,
    0.93188476858967106,
    0.63261795801119336,
    0.71603325282104247,
    0.41877472951438482,
    0.45029451839374096,
    0.55159041480602233,
    0.21753602570423891,
    0.70816525444993794,
    0.26607848480272622,
    0.46411755984161318,
    0.42320062138521458,
    0.37204549830118205,
    0.38968218796911197,
    0.11838791042231629,
    0.23460144710516165,
    0.34411427128675631,
    0.26618944631227137,
    0.51063751364053522,
    0.7793694420033549,
    0.59031136082192182,
    0.49195515928793193,
    0.40614376085860787,
    0.28144891399238653,
    0.43634815451011644,
    0.4788471911106883,
    0.53623690774424926,
    0.65717598975673075,
    0.63282415393669345,
    0.61792245577754978,
    0.68390616001298077,
    0.75938427746658632,
    0.7329978821339167,
    0.54448793302708565,
    0.62269387219607408,
    0.51084871034225199,
    0.50755833729248578,
    0.50476852743873268,
    0.40218982826403404,
    0.38992199911065317,
    0.5122740680925727,
    0.49535974919221276,
    0.41574240562373316,
    0.44898853295496383,
    0.45603919224192818,
    0.55085652015211994,
    0.54411015823651876,
    0.57396182234821236,
    0.6093070899113896,
    0.65280971383297948,
    0.53256440202563649,
    0.49913475680082209
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
